$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at the top of the "Vega Monumental Concepción" block (rows 376-379),
# shifting the existing rows 376-397 down to 380-401.
$ws.Rows("376:379").Insert()

# New weekly price data (week of 2022-02-18, serial 44610) for the 4 quality grades.
$newRows = @(
    @{ L = "1a amarillo"; M = 270; N = 17000; O = 18000; P = 17556; S = 1097 },
    @{ L = "1a plateado"; M = 220; N = 17000; O = 18000; P = 17545; S = 1097 },
    @{ L = "2a amarillo"; M = 220; N = 14000; O = 15000; P = 14455; S = 903 },
    @{ L = "2a plateado"; M = 180; N = 14000; O = 15000; P = 14556; S = 910 }
)

$r = 376
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = 11
    $ws.Cells.Item($r, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($r, 3).Value = "Bíobío"
    $ws.Cells.Item($r, 4).Value = 44610
    $ws.Cells.Item($r, 5).Value = 8
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100102
    $ws.Cells.Item($r, 8).Value = "Cítricos"
    $ws.Cells.Item($r, 9).Value = 100102003
    $ws.Cells.Item($r, 10).Value = "Limón"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = "`$/malla 16 kilos"
    $ws.Cells.Item($r, 18).Value = "Provincia de Melipilla"
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = 16
    $r = $r + 1
}
